$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct status name: "bleu" -> "noir"
$ws.Cells.Replace("bleu", "noir")

# Correct status wording (longer/more specific phrases first so they are not
# partially clobbered by the shorter base phrase replacement below)
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois")
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois")
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés")
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés")
